$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.328.03"
$ws.Range("E2").Value = "  +0.70%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.932.22"
$ws.Range("E3").Value = "  +0.53%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.52"
$ws.Range("E5").Value = "  +0.48%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.54"
$ws.Range("E6").Value = "  -1.54%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.31%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.20%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  +1.11%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -1.69%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -0.57%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -0.92%  "

# Row 13 - Avalanche
$ws.Range("E13").Value = "  -1.00%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.42%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.418.17"
$ws.Range("E15").Value = "  +0.57%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "61.320.59"
$ws.Range("E16").Value = "  +0.70%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.932.41"
$ws.Range("E17").Value = "  +0.50%  "

# Row 18 - Polkadot
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.63"
$ws.Range("E18").Value = "  -0.82%  "

# Row 19 - BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "433.65"
$ws.Range("E19").Value = "  +0.79%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.57"
$ws.Range("E20").Value = "  +1.53%  "

# Row 21 - Polygon
$ws.Range("E21").Value = "  -1.08%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.06"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.43"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24 - RenderToken
$ws.Range("E24").Value = "  -0.64%  "

# Row 25 - Fetch.AI
$ws.Range("E25").Value = "  -1.92%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.73"
$ws.Range("E26").Value = "  -1.78%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.01%  "

# Row 28 - ImmutableX
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -3.76%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -0.74%  "

# Row 30 - NEARProtocol
$ws.Range("E30").Value = "  -2.40%  "

# Row 31 - EthereumClassic
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.67"
$ws.Range("E31").Value = "  +0.90%  "

# Row 32 - Hedera
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.109"
$ws.Range("E32").Value = "  +1.32%  "

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  -0.01%  "

# Row 34 - PEPE
$ws.Range("D34").Value = "0.0₃0875"
$ws.Range("E34").Value = "  +2.74%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -0.27%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  -0.37%  "

# Row 37 - dogwifhat
$ws.Range("E37").Value = "  -1.89%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -0.43%  "

# Row 40 - Cosmos
$ws.Range("E40").Value = "  -0.69%  "

# Row 41 - Arweave
$ws.Range("E41").Value = "  +5.35%  "

# Row 42 - TheGraph
$ws.Range("E42").Value = "  -2.83%  "

# Row 43 - Maker
$ws.Range("D43").Value = "2.705.24"
$ws.Range("E43").Value = "  +0.12%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -0.32%  "

# Row 45 - (Monero->Bittensor)
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "364.59"
$ws.Range("E45").Value = "  -2.79%  "

# Row 46 - (Bittensor->Monero)
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.18"
$ws.Range("E46").Value = "  +0.94%  "

# Row 47 - USDe
$ws.Range("E47").Value = "  +0.06%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.56"
$ws.Range("E48").Value = "  -1.27%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -1.38%  "

# Row 50 - ThetaToken
$ws.Range("E50").Value = "  -0.45%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  +0.12%  "
